$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 121
$ws.Cells.Item(121, 2).Value = 6941361
$ws.Cells.Item(121, 6).Value = "Al Taee"
$ws.Cells.Item(121, 7).Value = "Al Raed"
$ws.Cells.Item(121, 8).Value = 4
$ws.Cells.Item(121, 9).Value = 3
$ws.Cells.Item(121, 10).Value = "H"
$ws.Cells.Item(121, 11).Value = 2.25
$ws.Cells.Item(121, 12).Value = 3.4
$ws.Cells.Item(121, 13).Value = 2.8
$ws.Cells.Item(121, 14).Value = 3.8
$ws.Cells.Item(121, 15).Value = 3.4
$ws.Cells.Item(121, 16).Value = 1.833
$ws.Cells.Item(121, 17).Value = 0.5
$ws.Cells.Item(121, 18).Value = 1.975
$ws.Cells.Item(121, 19).Value = 1.825
$ws.Cells.Item(121, 20).Value = 2.25
$ws.Cells.Item(121, 21).Value = 1.8
$ws.Cells.Item(121, 22).Value = 2
$ws.Cells.Item(121, 23).Value = 2.8
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(121, 25).Value = -1
$ws.Cells.Item(121, 26).Value = 0.9750000000000001
$ws.Cells.Item(121, 27).Value = -1
$ws.Cells.Item(121, 28).Value = 0.8
$ws.Cells.Item(121, 29).Value = -1

# Row 122
$ws.Cells.Item(122, 2).Value = 6940756
$ws.Cells.Item(122, 6).Value = "AlNassr Riyadh"
$ws.Cells.Item(122, 7).Value = "Al Akhdoud"
$ws.Cells.Item(122, 8).Value = 3
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = "H"
$ws.Cells.Item(122, 11).Value = 1.142
$ws.Cells.Item(122, 12).Value = 7
$ws.Cells.Item(122, 13).Value = 15
$ws.Cells.Item(122, 14).Value = 1.285
$ws.Cells.Item(122, 15).Value = 5.25
$ws.Cells.Item(122, 16).Value = 8.5
$ws.Cells.Item(122, 17).Value = -1.75
$ws.Cells.Item(122, 18).Value = 1.875
$ws.Cells.Item(122, 19).Value = 1.925
$ws.Cells.Item(122, 20).Value = 3.5
$ws.Cells.Item(122, 21).Value = 1.875
$ws.Cells.Item(122, 22).Value = 1.925
$ws.Cells.Item(122, 23).Value = 0.2849999999999999
$ws.Cells.Item(122, 24).Value = -1
$ws.Cells.Item(122, 25).Value = -1
$ws.Cells.Item(122, 26).Value = 0.875
$ws.Cells.Item(122, 27).Value = -1
$ws.Cells.Item(122, 28).Value = -1
$ws.Cells.Item(122, 29).Value = 0.925

# Row 125
$ws.Cells.Item(125, 2).Value = 6941362
$ws.Cells.Item(125, 6).Value = "Damac FC"
$ws.Cells.Item(125, 7).Value = "Abha"
$ws.Cells.Item(125, 8).Value = 4
$ws.Cells.Item(125, 9).Value = 2
$ws.Cells.Item(125, 10).Value = "H"
$ws.Cells.Item(125, 11).Value = 1.615
$ws.Cells.Item(125, 12).Value = 4
$ws.Cells.Item(125, 13).Value = 4.5
$ws.Cells.Item(125, 14).Value = 1.85
$ws.Cells.Item(125, 15).Value = 3.6
$ws.Cells.Item(125, 16).Value = 3.6
$ws.Cells.Item(125, 17).Value = -0.5
$ws.Cells.Item(125, 18).Value = 1.875
$ws.Cells.Item(125, 19).Value = 1.925
$ws.Cells.Item(125, 20).Value = 2.75
$ws.Cells.Item(125, 21).Value = 1.875
$ws.Cells.Item(125, 22).Value = 1.925
$ws.Cells.Item(125, 23).Value = 0.8500000000000001
$ws.Cells.Item(125, 24).Value = -1
$ws.Cells.Item(125, 25).Value = -1
$ws.Cells.Item(125, 26).Value = 0.875
$ws.Cells.Item(125, 27).Value = -1
$ws.Cells.Item(125, 28).Value = 0.875
$ws.Cells.Item(125, 29).Value = -1

# Row 126
$ws.Cells.Item(126, 2).Value = 6941360
$ws.Cells.Item(126, 6).Value = "Al Wehda Mecca"
$ws.Cells.Item(126, 7).Value = "Al Khaleej Saihat"
$ws.Cells.Item(126, 8).Value = 3
$ws.Cells.Item(126, 9).Value = 1
$ws.Cells.Item(126, 10).Value = "H"
$ws.Cells.Item(126, 11).Value = 1.8
$ws.Cells.Item(126, 12).Value = 3.6
$ws.Cells.Item(126, 13).Value = 3.9
$ws.Cells.Item(126, 14).Value = 2
$ws.Cells.Item(126, 15).Value = 3.4
$ws.Cells.Item(126, 16).Value = 3.4
$ws.Cells.Item(126, 17).Value = -0.5
$ws.Cells.Item(126, 18).Value = 2.025
$ws.Cells.Item(126, 19).Value = 1.775
$ws.Cells.Item(126, 20).Value = 2.75
$ws.Cells.Item(126, 21).Value = 2
$ws.Cells.Item(126, 22).Value = 1.8
$ws.Cells.Item(126, 23).Value = 1
$ws.Cells.Item(126, 24).Value = -1
$ws.Cells.Item(126, 25).Value = -1
$ws.Cells.Item(126, 26).Value = 1.025
$ws.Cells.Item(126, 27).Value = -1
$ws.Cells.Item(126, 28).Value = 1
$ws.Cells.Item(126, 29).Value = -1

# Row 128
$ws.Cells.Item(128, 2).Value = 6941374
$ws.Cells.Item(128, 6).Value = "Abha"
$ws.Cells.Item(128, 7).Value = "Al Ahli Jeddah"
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 6
$ws.Cells.Item(128, 10).Value = "A"
$ws.Cells.Item(128, 11).Value = 4.5
$ws.Cells.Item(128, 12).Value = 4
$ws.Cells.Item(128, 13).Value = 1.615
$ws.Cells.Item(128, 14).Value = 4.75
$ws.Cells.Item(128, 15).Value = 4.2
$ws.Cells.Item(128, 16).Value = 1.571
$ws.Cells.Item(128, 17).Value = 1
$ws.Cells.Item(128, 18).Value = 1.925
$ws.Cells.Item(128, 19).Value = 1.875
$ws.Cells.Item(128, 20).Value = 3.25
$ws.Cells.Item(128, 21).Value = 1.95
$ws.Cells.Item(128, 22).Value = 1.85
$ws.Cells.Item(128, 23).Value = -1
$ws.Cells.Item(128, 24).Value = -1
$ws.Cells.Item(128, 25).Value = 0.571
$ws.Cells.Item(128, 26).Value = -1
$ws.Cells.Item(128, 27).Value = 0.875
$ws.Cells.Item(128, 28).Value = 0.95
$ws.Cells.Item(128, 29).Value = -1

# Row 129
$ws.Cells.Item(129, 2).Value = 6941371
$ws.Cells.Item(129, 6).Value = "Al Riyadh"
$ws.Cells.Item(129, 7).Value = "Al Hazm"
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = "D"
$ws.Cells.Item(129, 11).Value = 1.75
$ws.Cells.Item(129, 12).Value = 3.75
$ws.Cells.Item(129, 13).Value = 4
$ws.Cells.Item(129, 14).Value = 2.7
$ws.Cells.Item(129, 15).Value = 3.2
$ws.Cells.Item(129, 16).Value = 2.5
$ws.Cells.Item(129, 17).Value = 0
$ws.Cells.Item(129, 18).Value = 2
$ws.Cells.Item(129, 19).Value = 1.8
$ws.Cells.Item(129, 20).Value = 2.5
$ws.Cells.Item(129, 21).Value = 1.875
$ws.Cells.Item(129, 22).Value = 1.925
$ws.Cells.Item(129, 23).Value = -1
$ws.Cells.Item(129, 24).Value = 2.2
$ws.Cells.Item(129, 25).Value = -1
$ws.Cells.Item(129, 26).Value = 0
$ws.Cells.Item(129, 27).Value = 0
$ws.Cells.Item(129, 28).Value = -1
$ws.Cells.Item(129, 29).Value = 0.925

# Row 137
$ws.Cells.Item(137, 2).Value = 6940761
$ws.Cells.Item(137, 6).Value = "Damac FC"
$ws.Cells.Item(137, 7).Value = "Al Ittihad Jeddah"
$ws.Cells.Item(137, 8).Value = 3
$ws.Cells.Item(137, 9).Value = 1
$ws.Cells.Item(137, 10).Value = "H"
$ws.Cells.Item(137, 11).Value = 4.5
$ws.Cells.Item(137, 12).Value = 4.2
$ws.Cells.Item(137, 13).Value = 1.6
$ws.Cells.Item(137, 14).Value = 2.9
$ws.Cells.Item(137, 15).Value = 3.5
$ws.Cells.Item(137, 16).Value = 2.15
$ws.Cells.Item(137, 17).Value = 0.25
$ws.Cells.Item(137, 18).Value = 1.825
$ws.Cells.Item(137, 19).Value = 1.975
$ws.Cells.Item(137, 20).Value = 2.75
$ws.Cells.Item(137, 21).Value = 2
$ws.Cells.Item(137, 22).Value = 1.8
$ws.Cells.Item(137, 23).Value = 1.9
$ws.Cells.Item(137, 24).Value = -1
$ws.Cells.Item(137, 25).Value = -1
$ws.Cells.Item(137, 26).Value = 0.825
$ws.Cells.Item(137, 27).Value = -1
$ws.Cells.Item(137, 28).Value = 1
$ws.Cells.Item(137, 29).Value = -1

# Row 138
$ws.Cells.Item(138, 2).Value = 6941375
$ws.Cells.Item(138, 6).Value = "Al Hazm"
$ws.Cells.Item(138, 7).Value = "Al Fateh SC"
$ws.Cells.Item(138, 8).Value = 2
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = "H"
$ws.Cells.Item(138, 11).Value = 5
$ws.Cells.Item(138, 12).Value = 3.9
$ws.Cells.Item(138, 13).Value = 1.6
$ws.Cells.Item(138, 14).Value = 3.1
$ws.Cells.Item(138, 15).Value = 3.5
$ws.Cells.Item(138, 16).Value = 2.15
$ws.Cells.Item(138, 17).Value = 0.25
$ws.Cells.Item(138, 18).Value = 1.925
$ws.Cells.Item(138, 19).Value = 1.875
$ws.Cells.Item(138, 20).Value = 2.75
$ws.Cells.Item(138, 21).Value = 1.825
$ws.Cells.Item(138, 22).Value = 1.975
$ws.Cells.Item(138, 23).Value = 2.1
$ws.Cells.Item(138, 24).Value = -1
$ws.Cells.Item(138, 25).Value = -1
$ws.Cells.Item(138, 26).Value = 0.925
$ws.Cells.Item(138, 27).Value = -1
$ws.Cells.Item(138, 28).Value = -1
$ws.Cells.Item(138, 29).Value = 0.9750000000000001

# Row 140
$ws.Cells.Item(140, 2).Value = 6940759
$ws.Cells.Item(140, 6).Value = "Al Taee"
$ws.Cells.Item(140, 7).Value = "Al Hilal Riyadh"
$ws.Cells.Item(140, 8).Value = 1
$ws.Cells.Item(140, 9).Value = 2
$ws.Cells.Item(140, 10).Value = "A"
$ws.Cells.Item(140, 11).Value = 9.5
$ws.Cells.Item(140, 12).Value = 6.5
$ws.Cells.Item(140, 13).Value = 1.2
$ws.Cells.Item(140, 14).Value = 15
$ws.Cells.Item(140, 15).Value = 9
$ws.Cells.Item(140, 16).Value = 1.111
$ws.Cells.Item(140, 17).Value = 2.5
$ws.Cells.Item(140, 18).Value = 1.9
$ws.Cells.Item(140, 19).Value = 1.9
$ws.Cells.Item(140, 20).Value = 3.75
$ws.Cells.Item(140, 21).Value = 1.975
$ws.Cells.Item(140, 22).Value = 1.825
$ws.Cells.Item(140, 23).Value = -1
$ws.Cells.Item(140, 24).Value = -1
$ws.Cells.Item(140, 25).Value = 0.111
$ws.Cells.Item(140, 26).Value = 0.8999999999999999
$ws.Cells.Item(140, 27).Value = -1
$ws.Cells.Item(140, 28).Value = -1
$ws.Cells.Item(140, 29).Value = 0.825

# Row 142
$ws.Cells.Item(142, 2).Value = 6941379
$ws.Cells.Item(142, 6).Value = "Al Khaleej Saihat"
$ws.Cells.Item(142, 7).Value = "Abha"
$ws.Cells.Item(142, 8).Value = 3
$ws.Cells.Item(142, 9).Value = 1
$ws.Cells.Item(142, 10).Value = "H"
$ws.Cells.Item(142, 11).Value = 2.15
$ws.Cells.Item(142, 12).Value = 3.5
$ws.Cells.Item(142, 13).Value = 3
$ws.Cells.Item(142, 14).Value = 1.75
$ws.Cells.Item(142, 15).Value = 3.8
$ws.Cells.Item(142, 16).Value = 4
$ws.Cells.Item(142, 17).Value = -0.75
$ws.Cells.Item(142, 18).Value = 1.975
$ws.Cells.Item(142, 19).Value = 1.825
$ws.Cells.Item(142, 20).Value = 3
$ws.Cells.Item(142, 21).Value = 2
$ws.Cells.Item(142, 22).Value = 1.8
$ws.Cells.Item(142, 23).Value = 0.75
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = -1
$ws.Cells.Item(142, 26).Value = 0.9750000000000001
$ws.Cells.Item(142, 27).Value = -1
$ws.Cells.Item(142, 28).Value = 1
$ws.Cells.Item(142, 29).Value = -1

# Row 150
$ws.Cells.Item(150, 2).Value = 6940763
$ws.Cells.Item(150, 6).Value = "Al Hilal Riyadh"
$ws.Cells.Item(150, 7).Value = "Al Wehda Mecca"
$ws.Cells.Item(150, 8).Value = 2
$ws.Cells.Item(150, 9).Value = 0
$ws.Cells.Item(150, 10).Value = "H"
$ws.Cells.Item(150, 11).Value = 1.2
$ws.Cells.Item(150, 12).Value = 7
$ws.Cells.Item(150, 13).Value = 9.5
$ws.Cells.Item(150, 14).Value = 1.222
$ws.Cells.Item(150, 15).Value = 7
$ws.Cells.Item(150, 16).Value = 8.5
$ws.Cells.Item(150, 17).Value = -2
$ws.Cells.Item(150, 18).Value = 1.85
$ws.Cells.Item(150, 19).Value = 1.95
$ws.Cells.Item(150, 20).Value = 3.75
$ws.Cells.Item(150, 21).Value = 1.9
$ws.Cells.Item(150, 22).Value = 1.9
$ws.Cells.Item(150, 23).Value = 0.222
$ws.Cells.Item(150, 24).Value = -1
$ws.Cells.Item(150, 25).Value = -1
$ws.Cells.Item(150, 26).Value = 0
$ws.Cells.Item(150, 27).Value = 0
$ws.Cells.Item(150, 28).Value = -1
$ws.Cells.Item(150, 29).Value = 0.8999999999999999

# Row 151
$ws.Cells.Item(151, 2).Value = 6940818
$ws.Cells.Item(151, 6).Value = "Al Fateh SC"
$ws.Cells.Item(151, 7).Value = "Al Shabab Riyadh"
$ws.Cells.Item(151, 8).Value = 1
$ws.Cells.Item(151, 9).Value = 1
$ws.Cells.Item(151, 10).Value = "D"
$ws.Cells.Item(151, 11).Value = 2.55
$ws.Cells.Item(151, 12).Value = 3.4
$ws.Cells.Item(151, 13).Value = 2.45
$ws.Cells.Item(151, 14).Value = 2.375
$ws.Cells.Item(151, 15).Value = 3.5
$ws.Cells.Item(151, 16).Value = 2.625
$ws.Cells.Item(151, 17).Value = 0
$ws.Cells.Item(151, 18).Value = 1.8
$ws.Cells.Item(151, 19).Value = 2
$ws.Cells.Item(151, 20).Value = 2.75
$ws.Cells.Item(151, 21).Value = 1.775
$ws.Cells.Item(151, 22).Value = 2.025
$ws.Cells.Item(151, 23).Value = -1
$ws.Cells.Item(151, 24).Value = 2.5
$ws.Cells.Item(151, 25).Value = -1
$ws.Cells.Item(151, 26).Value = 0
$ws.Cells.Item(151, 27).Value = 0
$ws.Cells.Item(151, 28).Value = -1
$ws.Cells.Item(151, 29).Value = 1.025

# Row 157
$ws.Cells.Item(157, 2).Value = 6940819
$ws.Cells.Item(157, 6).Value = "Al Akhdoud"
$ws.Cells.Item(157, 7).Value = "Al Shabab Riyadh"
$ws.Cells.Item(157, 8).Value = 1
$ws.Cells.Item(157, 9).Value = 0
$ws.Cells.Item(157, 10).Value = "H"
$ws.Cells.Item(157, 11).Value = 2.5
$ws.Cells.Item(157, 12).Value = 3.4
$ws.Cells.Item(157, 13).Value = 2.6
$ws.Cells.Item(157, 14).Value = 2.7
$ws.Cells.Item(157, 15).Value = 3.25
$ws.Cells.Item(157, 16).Value = 2.5
$ws.Cells.Item(157, 17).Value = 0
$ws.Cells.Item(157, 18).Value = 1.95
$ws.Cells.Item(157, 19).Value = 1.85
$ws.Cells.Item(157, 20).Value = 2.5
$ws.Cells.Item(157, 21).Value = 2
$ws.Cells.Item(157, 22).Value = 1.8
$ws.Cells.Item(157, 23).Value = 1.7
$ws.Cells.Item(157, 24).Value = -1
$ws.Cells.Item(157, 25).Value = -1
$ws.Cells.Item(157, 26).Value = 0.95
$ws.Cells.Item(157, 27).Value = -1
$ws.Cells.Item(157, 28).Value = -1
$ws.Cells.Item(157, 29).Value = 0.8

# Row 158
$ws.Cells.Item(158, 2).Value = 7605556
$ws.Cells.Item(158, 6).Value = "AlNassr Riyadh"
$ws.Cells.Item(158, 7).Value = "Al Ittifaq Dammam"
$ws.Cells.Item(158, 8).Value = 3
$ws.Cells.Item(158, 9).Value = 1
$ws.Cells.Item(158, 10).Value = "H"
$ws.Cells.Item(158, 11).Value = 1.285
$ws.Cells.Item(158, 12).Value = 5.25
$ws.Cells.Item(158, 13).Value = 7.5
$ws.Cells.Item(158, 14).Value = 1.333
$ws.Cells.Item(158, 15).Value = 5
$ws.Cells.Item(158, 16).Value = 7
$ws.Cells.Item(158, 17).Value = -1.5
$ws.Cells.Item(158, 18).Value = 1.925
$ws.Cells.Item(158, 19).Value = 1.875
$ws.Cells.Item(158, 20).Value = 3
$ws.Cells.Item(158, 21).Value = 1.85
$ws.Cells.Item(158, 22).Value = 1.95
$ws.Cells.Item(158, 23).Value = 0.333
$ws.Cells.Item(158, 24).Value = -1
$ws.Cells.Item(158, 25).Value = -1
$ws.Cells.Item(158, 26).Value = 0.925
$ws.Cells.Item(158, 27).Value = -1
$ws.Cells.Item(158, 28).Value = 0.8500000000000001
$ws.Cells.Item(158, 29).Value = -1

# Row 164
$ws.Cells.Item(164, 2).Value = 6941393
$ws.Cells.Item(164, 6).Value = "Al Ittifaq Dammam"
$ws.Cells.Item(164, 7).Value = "Al Hazm"
$ws.Cells.Item(164, 8).Value = 1
$ws.Cells.Item(164, 9).Value = 1
$ws.Cells.Item(164, 10).Value = "D"
$ws.Cells.Item(164, 11).Value = 1.571
$ws.Cells.Item(164, 12).Value = 4
$ws.Cells.Item(164, 13).Value = 5
$ws.Cells.Item(164, 14).Value = 1.533
$ws.Cells.Item(164, 15).Value = 4
$ws.Cells.Item(164, 16).Value = 5.5
$ws.Cells.Item(164, 17).Value = -1
$ws.Cells.Item(164, 18).Value = 1.9
$ws.Cells.Item(164, 19).Value = 1.9
$ws.Cells.Item(164, 20).Value = 2.5
$ws.Cells.Item(164, 21).Value = 1.8
$ws.Cells.Item(164, 22).Value = 2
$ws.Cells.Item(164, 23).Value = -1
$ws.Cells.Item(164, 24).Value = 3
$ws.Cells.Item(164, 25).Value = -1
$ws.Cells.Item(164, 26).Value = -1
$ws.Cells.Item(164, 27).Value = 0.8999999999999999
$ws.Cells.Item(164, 28).Value = -1
$ws.Cells.Item(164, 29).Value = 1

# Row 165
$ws.Cells.Item(165, 2).Value = 6941395
$ws.Cells.Item(165, 6).Value = "Al Riyadh"
$ws.Cells.Item(165, 7).Value = "Damac FC"
$ws.Cells.Item(165, 8).Value = 1
$ws.Cells.Item(165, 9).Value = 0
$ws.Cells.Item(165, 10).Value = "H"
$ws.Cells.Item(165, 11).Value = 3.75
$ws.Cells.Item(165, 12).Value = 3.6
$ws.Cells.Item(165, 13).Value = 1.833
$ws.Cells.Item(165, 14).Value = 4.333
$ws.Cells.Item(165, 15).Value = 3.4
$ws.Cells.Item(165, 16).Value = 1.75
$ws.Cells.Item(165, 17).Value = 0.5
$ws.Cells.Item(165, 18).Value = 2
$ws.Cells.Item(165, 19).Value = 1.8
$ws.Cells.Item(165, 20).Value = 2.5
$ws.Cells.Item(165, 21).Value = 1.95
$ws.Cells.Item(165, 22).Value = 1.85
$ws.Cells.Item(165, 23).Value = 3.333
$ws.Cells.Item(165, 24).Value = -1
$ws.Cells.Item(165, 25).Value = -1
$ws.Cells.Item(165, 26).Value = 1
$ws.Cells.Item(165, 27).Value = -1
$ws.Cells.Item(165, 28).Value = -1
$ws.Cells.Item(165, 29).Value = 0.8500000000000001

# Row 191
$ws.Cells.Item(191, 2).Value = 6941407
$ws.Cells.Item(191, 6).Value = "Al Wehda Mecca"
$ws.Cells.Item(191, 7).Value = "Al Taawon Buraidah"
$ws.Cells.Item(191, 8).Value = 3
$ws.Cells.Item(191, 9).Value = 3
$ws.Cells.Item(191, 10).Value = "D"
$ws.Cells.Item(191, 11).Value = 2.2
$ws.Cells.Item(191, 12).Value = 3.5
$ws.Cells.Item(191, 13).Value = 3.1
$ws.Cells.Item(191, 14).Value = 2.8
$ws.Cells.Item(191, 15).Value = 3.5
$ws.Cells.Item(191, 16).Value = 2.375
$ws.Cells.Item(191, 17).Value = 0.25
$ws.Cells.Item(191, 18).Value = 1.75
$ws.Cells.Item(191, 19).Value = 2.05
$ws.Cells.Item(191, 20).Value = 2.75
$ws.Cells.Item(191, 21).Value = 1.925
$ws.Cells.Item(191, 22).Value = 1.875
$ws.Cells.Item(191, 23).Value = -1
$ws.Cells.Item(191, 24).Value = 2.5
$ws.Cells.Item(191, 25).Value = -1
$ws.Cells.Item(191, 26).Value = 0.375
$ws.Cells.Item(191, 27).Value = -0.5
$ws.Cells.Item(191, 28).Value = 0.925
$ws.Cells.Item(191, 29).Value = -1

# Row 192
$ws.Cells.Item(192, 2).Value = 6941409
$ws.Cells.Item(192, 6).Value = "Al Raed"
$ws.Cells.Item(192, 7).Value = "Al Fayha"
$ws.Cells.Item(192, 8).Value = 1
$ws.Cells.Item(192, 9).Value = 2
$ws.Cells.Item(192, 10).Value = "A"
$ws.Cells.Item(192, 11).Value = 2.05
$ws.Cells.Item(192, 12).Value = 3.4
$ws.Cells.Item(192, 13).Value = 3.5
$ws.Cells.Item(192, 14).Value = 2.2
$ws.Cells.Item(192, 15).Value = 3.3
$ws.Cells.Item(192, 16).Value = 3.25
$ws.Cells.Item(192, 17).Value = -0.25
$ws.Cells.Item(192, 18).Value = 1.9
$ws.Cells.Item(192, 19).Value = 1.9
$ws.Cells.Item(192, 20).Value = 2.25
$ws.Cells.Item(192, 21).Value = 1.825
$ws.Cells.Item(192, 22).Value = 1.975
$ws.Cells.Item(192, 23).Value = -1
$ws.Cells.Item(192, 24).Value = -1
$ws.Cells.Item(192, 25).Value = 2.25
$ws.Cells.Item(192, 26).Value = -1
$ws.Cells.Item(192, 27).Value = 0.8999999999999999
$ws.Cells.Item(192, 28).Value = 0.825
$ws.Cells.Item(192, 29).Value = -1

# Row 194
$ws.Cells.Item(194, 2).Value = 6941410
$ws.Cells.Item(194, 6).Value = "Al Akhdoud"
$ws.Cells.Item(194, 7).Value = "Al Riyadh"
$ws.Cells.Item(194, 8).Value = 1
$ws.Cells.Item(194, 9).Value = 2
$ws.Cells.Item(194, 10).Value = "A"
$ws.Cells.Item(194, 11).Value = 1.615
$ws.Cells.Item(194, 12).Value = 4
$ws.Cells.Item(194, 13).Value = 5
$ws.Cells.Item(194, 14).Value = 1.8
$ws.Cells.Item(194, 15).Value = 3.5
$ws.Cells.Item(194, 16).Value = 4.5
$ws.Cells.Item(194, 17).Value = -0.5
$ws.Cells.Item(194, 18).Value = 1.775
$ws.Cells.Item(194, 19).Value = 2.025
$ws.Cells.Item(194, 20).Value = 2.25
$ws.Cells.Item(194, 21).Value = 1.9
$ws.Cells.Item(194, 22).Value = 1.9
$ws.Cells.Item(194, 23).Value = -1
$ws.Cells.Item(194, 24).Value = -1
$ws.Cells.Item(194, 25).Value = 3.5
$ws.Cells.Item(194, 26).Value = -1
$ws.Cells.Item(194, 27).Value = 1.025
$ws.Cells.Item(194, 28).Value = 0.8999999999999999
$ws.Cells.Item(194, 29).Value = -1

# Row 195
$ws.Cells.Item(195, 2).Value = 6941412
$ws.Cells.Item(195, 6).Value = "Abha"
$ws.Cells.Item(195, 7).Value = "Al Taee"
$ws.Cells.Item(195, 8).Value = 2
$ws.Cells.Item(195, 9).Value = 0
$ws.Cells.Item(195, 10).Value = "H"
$ws.Cells.Item(195, 11).Value = 2.4
$ws.Cells.Item(195, 12).Value = 3.5
$ws.Cells.Item(195, 13).Value = 2.75
$ws.Cells.Item(195, 14).Value = 2.2
$ws.Cells.Item(195, 15).Value = 3.6
$ws.Cells.Item(195, 16).Value = 3
$ws.Cells.Item(195, 17).Value = -0.25
$ws.Cells.Item(195, 18).Value = 1.975
$ws.Cells.Item(195, 19).Value = 1.825
$ws.Cells.Item(195, 20).Value = 2.75
$ws.Cells.Item(195, 21).Value = 1.925
$ws.Cells.Item(195, 22).Value = 1.875
$ws.Cells.Item(195, 23).Value = 1.2
$ws.Cells.Item(195, 24).Value = -1
$ws.Cells.Item(195, 25).Value = -1
$ws.Cells.Item(195, 26).Value = 0.9750000000000001
$ws.Cells.Item(195, 27).Value = -1
$ws.Cells.Item(195, 28).Value = -1
$ws.Cells.Item(195, 29).Value = 0.875

# Row 210
$ws.Cells.Item(210, 2).Value = 6940824
$ws.Cells.Item(210, 6).Value = "Al Hazm"
$ws.Cells.Item(210, 7).Value = "Al Shabab Riyadh"
$ws.Cells.Item(210, 8).Value = 0
$ws.Cells.Item(210, 9).Value = 3
$ws.Cells.Item(210, 10).Value = "A"
$ws.Cells.Item(210, 11).Value = 4.75
$ws.Cells.Item(210, 12).Value = 3.75
$ws.Cells.Item(210, 13).Value = 1.666
$ws.Cells.Item(210, 14).Value = 4.2
$ws.Cells.Item(210, 15).Value = 3.6
$ws.Cells.Item(210, 16).Value = 1.8
$ws.Cells.Item(210, 17).Value = 0.5
$ws.Cells.Item(210, 18).Value = 2
$ws.Cells.Item(210, 19).Value = 1.8
$ws.Cells.Item(210, 20).Value = 2.75
$ws.Cells.Item(210, 21).Value = 1.95
$ws.Cells.Item(210, 22).Value = 1.85
$ws.Cells.Item(210, 23).Value = -1
$ws.Cells.Item(210, 24).Value = -1
$ws.Cells.Item(210, 25).Value = 0.8
$ws.Cells.Item(210, 26).Value = -1
$ws.Cells.Item(210, 27).Value = 0.8
$ws.Cells.Item(210, 28).Value = 0.475
$ws.Cells.Item(210, 29).Value = -0.5

# Row 211
$ws.Cells.Item(211, 2).Value = 6941418
$ws.Cells.Item(211, 6).Value = "Al Khaleej Saihat"
$ws.Cells.Item(211, 7).Value = "Al Akhdoud"
$ws.Cells.Item(211, 8).Value = 2
$ws.Cells.Item(211, 9).Value = 2
$ws.Cells.Item(211, 10).Value = "D"
$ws.Cells.Item(211, 11).Value = 2.1
$ws.Cells.Item(211, 12).Value = 3.2
$ws.Cells.Item(211, 13).Value = 3.4
$ws.Cells.Item(211, 14).Value = 2
$ws.Cells.Item(211, 15).Value = 3.3
$ws.Cells.Item(211, 16).Value = 3.6
$ws.Cells.Item(211, 17).Value = -0.5
$ws.Cells.Item(211, 18).Value = 2.025
$ws.Cells.Item(211, 19).Value = 1.775
$ws.Cells.Item(211, 20).Value = 2.5
$ws.Cells.Item(211, 21).Value = 1.925
$ws.Cells.Item(211, 22).Value = 1.875
$ws.Cells.Item(211, 23).Value = -1
$ws.Cells.Item(211, 24).Value = 2.3
$ws.Cells.Item(211, 25).Value = -1
$ws.Cells.Item(211, 26).Value = -1
$ws.Cells.Item(211, 27).Value = 0.7749999999999999
$ws.Cells.Item(211, 28).Value = 0.925
$ws.Cells.Item(211, 29).Value = -1

# Row 216
$ws.Cells.Item(216, 2).Value = 6940783
$ws.Cells.Item(216, 6).Value = "Al Fateh SC"
$ws.Cells.Item(216, 7).Value = "Al Ittihad Jeddah"
$ws.Cells.Item(216, 8).Value = 2
$ws.Cells.Item(216, 9).Value = 4
$ws.Cells.Item(216, 10).Value = "A"
$ws.Cells.Item(216, 11).Value = 3.75
$ws.Cells.Item(216, 12).Value = 3.75
$ws.Cells.Item(216, 13).Value = 1.85
$ws.Cells.Item(216, 14).Value = 4.2
$ws.Cells.Item(216, 15).Value = 3.6
$ws.Cells.Item(216, 16).Value = 1.8
$ws.Cells.Item(216, 17).Value = 0.5
$ws.Cells.Item(216, 18).Value = 1.975
$ws.Cells.Item(216, 19).Value = 1.825
$ws.Cells.Item(216, 20).Value = 2.75
$ws.Cells.Item(216, 21).Value = 1.9
$ws.Cells.Item(216, 22).Value = 1.9
$ws.Cells.Item(216, 23).Value = -1
$ws.Cells.Item(216, 24).Value = -1
$ws.Cells.Item(216, 25).Value = 0.8
$ws.Cells.Item(216, 26).Value = -1
$ws.Cells.Item(216, 27).Value = 0.825
$ws.Cells.Item(216, 28).Value = 0.8999999999999999
$ws.Cells.Item(216, 29).Value = -1

# Row 217
$ws.Cells.Item(217, 2).Value = 6941419
$ws.Cells.Item(217, 6).Value = "Al Raed"
$ws.Cells.Item(217, 7).Value = "Al Taawon Buraidah"
$ws.Cells.Item(217, 8).Value = 0
$ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(217, 10).Value = "D"
$ws.Cells.Item(217, 11).Value = 2.9
$ws.Cells.Item(217, 12).Value = 3.4
$ws.Cells.Item(217, 13).Value = 2.3
$ws.Cells.Item(217, 14).Value = 2.55
$ws.Cells.Item(217, 15).Value = 3.3
$ws.Cells.Item(217, 16).Value = 2.7
$ws.Cells.Item(217, 17).Value = 0
$ws.Cells.Item(217, 18).Value = 1.85
$ws.Cells.Item(217, 19).Value = 1.95
$ws.Cells.Item(217, 20).Value = 2.5
$ws.Cells.Item(217, 21).Value = 1.825
$ws.Cells.Item(217, 22).Value = 1.975
$ws.Cells.Item(217, 23).Value = -1
$ws.Cells.Item(217, 24).Value = 2.3
$ws.Cells.Item(217, 25).Value = -1
$ws.Cells.Item(217, 26).Value = 0
$ws.Cells.Item(217, 27).Value = 0
$ws.Cells.Item(217, 28).Value = -1
$ws.Cells.Item(217, 29).Value = 0.9750000000000001

# Row 234
$ws.Cells.Item(234, 2).Value = 6941441
$ws.Cells.Item(234, 6).Value = "Damac FC"
$ws.Cells.Item(234, 7).Value = "AlNassr Riyadh"
$ws.Cells.Item(234, 11).Value = 4.333
$ws.Cells.Item(234, 12).Value = 4.75
$ws.Cells.Item(234, 13).Value = 1.55
$ws.Cells.Item(234, 14).Value = 5.25
$ws.Cells.Item(234, 15).Value = 5
$ws.Cells.Item(234, 16).Value = 1.444
$ws.Cells.Item(234, 17).Value = 1.25
$ws.Cells.Item(234, 18).Value = 1.85
$ws.Cells.Item(234, 19).Value = 1.95
$ws.Cells.Item(234, 20).Value = 3.25
$ws.Cells.Item(234, 21).Value = 1.875
$ws.Cells.Item(234, 22).Value = 1.925

# Row 235
$ws.Cells.Item(235, 21).Value = 1.975
$ws.Cells.Item(235, 22).Value = 1.825

# Row 236
$ws.Cells.Item(236, 2).Value = 6940788
$ws.Cells.Item(236, 6).Value = "Al Ittihad Jeddah"
$ws.Cells.Item(236, 7).Value = "Al Taawon Buraidah"
$ws.Cells.Item(236, 11).Value = 1.727
$ws.Cells.Item(236, 12).Value = 3.6
$ws.Cells.Item(236, 13).Value = 4.333
$ws.Cells.Item(236, 14).Value = 1.727
$ws.Cells.Item(236, 15).Value = 3.6
$ws.Cells.Item(236, 16).Value = 4.333
$ws.Cells.Item(236, 17).Value = -0.75
$ws.Cells.Item(236, 18).Value = 1.95
$ws.Cells.Item(236, 19).Value = 1.85
$ws.Cells.Item(236, 20).Value = 2.75
$ws.Cells.Item(236, 21).Value = 1.875
$ws.Cells.Item(236, 22).Value = 1.925

# Row 237
$ws.Cells.Item(237, 14).Value = 8.5
$ws.Cells.Item(237, 15).Value = 5.25
$ws.Cells.Item(237, 16).Value = 1.285
$ws.Cells.Item(237, 17).Value = 1.5
$ws.Cells.Item(237, 20).Value = 3
$ws.Cells.Item(237, 21).Value = 1.9
$ws.Cells.Item(237, 22).Value = 1.9

# Row 238
$ws.Cells.Item(238, 16).Value = 7.5
$ws.Cells.Item(238, 18).Value = 2
$ws.Cells.Item(238, 19).Value = 1.8
$ws.Cells.Item(238, 21).Value = 2
$ws.Cells.Item(238, 22).Value = 1.8

# Row 239
$ws.Cells.Item(239, 14).Value = 4.333
$ws.Cells.Item(239, 15).Value = 3.8
$ws.Cells.Item(239, 16).Value = 1.65
$ws.Cells.Item(239, 18).Value = 1.975
$ws.Cells.Item(239, 19).Value = 1.825
$ws.Cells.Item(239, 21).Value = 1.825
$ws.Cells.Item(239, 22).Value = 1.975

# Row 240
$ws.Cells.Item(240, 14).Value = 3.75
$ws.Cells.Item(240, 18).Value = 1.95
$ws.Cells.Item(240, 19).Value = 1.85
$ws.Cells.Item(240, 20).Value = 3
$ws.Cells.Item(240, 21).Value = 2
$ws.Cells.Item(240, 22).Value = 1.8

# Row 241
$ws.Cells.Item(241, 14).Value = 1.615

# Row 242
$ws.Cells.Item(242, 14).Value = 2.375
$ws.Cells.Item(242, 15).Value = 3.25
$ws.Cells.Item(242, 16).Value = 2.75
$ws.Cells.Item(242, 18).Value = 1.75
$ws.Cells.Item(242, 19).Value = 2.05
